# Remove the first 5 data rows (Trial_No 103-107) from the schedule sheet.
# Everything below shifts up by 5 rows; Trial_No values are untouched because
# they already belonged to the rows that remain (108-215).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A6").EntireRow.Delete() | Out-Null

# Re-establish the AutoFilter over the new (shrunk) data extent so the
# worksheet's <autoFilter> ref and the workbook's hidden _FilterDatabase
# defined name both collapse from $A$1:$K$114 to $A$1:$K$109.
$ws.AutoFilterMode = $false
$ws.Range("A1:K109").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$K`$109"
    }
}

# Update the saved selection/active cell shown when the workbook is reopened.
$ws.Range("D5").Select() | Out-Null
